$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header F1: "pvalues" -> "p"
$ws.Range("F1").Value = "p"

# Updated p-values after allowing MR (raw and SPSS) corrections
# (leading apostrophe forces text storage so values stay ".401"/".519"
# instead of being coerced into numbers, matching the original inlineStr text)
$ws.Range("F2").Value = "'.401"
$ws.Range("F3").Value = "'.401"
$ws.Range("F4").Value = "'.401"
$ws.Range("F6").Value = "'.519"

# New footnote row describing the multiple-tests correction applied
$ws.Range("A9").Value = "Multiple tests correction applied to p values: Benjamini-Hochberg"
